# Korjattu virhe muuttujalistan alustuksessa
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix misspelled names (missing Finnish/Swedish umlauts)
$ws.Range("A1").Value  = "Arvi Syrjänen"
$ws.Range("A2").Value  = "Ella Seppä"
$ws.Range("A7").Value  = "Jari Leppänen"
$ws.Range("A20").Value = "Paula Jyrkönen"

# Widen column A so the names fit better
$ws.Columns("A").ColumnWidth = 26

# Row 20 no longer needs to wrap onto extra lines once the column is wider
$ws.Rows(20).RowHeight = 15

# Update the view: scroll position and the active selection
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("B23").Select()
